$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1666.6666
$ws.Range("I40").Value = 1500
$ws.Range("K40").Value = 1500
$ws.Range("M40").Value = -1325
$ws.Range("H64").Value = 3331
$ws.Range("I64").Value = 3222.8462
$ws.Range("J64").Value = 3405
$ws.Range("K64").Value = 3222.8462
$ws.Range("L64").Value = 3405
$ws.Range("M64").Value = -2974.8462
$ws.Range("N64").Value = -3901
$ws.Range("H67").Value = 3331
$ws.Range("I67").Value = 3222.8462
$ws.Range("J67").Value = 3405
$ws.Range("K67").Value = 3222.8462
$ws.Range("L67").Value = 3405
$ws.Range("M67").Value = -2364.8462
$ws.Range("N67").Value = -5121
$ws.Range("H74").Value = 4041.25
$ws.Range("I74").Value = 3497.5
$ws.Range("J74").Value = 4150
$ws.Range("K74").Value = 3497.5
$ws.Range("L74").Value = 4150
$ws.Range("M74").Value = -2561.5
$ws.Range("N74").Value = -6022
$ws.Range("H76").Value = 4463.846
$ws.Range("I76").Value = 4503.3335
$ws.Range("J76").Value = 4375
$ws.Range("K76").Value = 4503.3335
$ws.Range("L76").Value = 4375
$ws.Range("M76").Value = -4188.3335
$ws.Range("N76").Value = -5005
$ws.Range("H77").Value = 4041.25
$ws.Range("I77").Value = 3497.5
$ws.Range("J77").Value = 4150
$ws.Range("K77").Value = 17487.5
$ws.Range("L77").Value = 20750
$ws.Range("M77").Value = -12807.5
$ws.Range("N77").Value = -30110
$ws.Range("H79").Value = 4463.846
$ws.Range("I79").Value = 4503.3335
$ws.Range("J79").Value = 4375
$ws.Range("K79").Value = 4503.3335
$ws.Range("L79").Value = 4375
$ws.Range("M79").Value = -3411.3335
$ws.Range("N79").Value = -6559
$ws.Range("H138").Value = 1764.22
$ws.Range("I138").Value = 835.38464
$ws.Range("J138").Value = 2090.5676
$ws.Range("K138").Value = 2506.15392
$ws.Range("L138").Value = 6271.702799999999
$ws.Range("M138").Value = 2633.84608
$ws.Range("N138").Value = -16551.7028

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18874.523
$ws.Range("I32").Value = 18150.674
$ws.Range("K32").Value = 18150.674
$ws.Range("M32").Value = -17863.674
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = $null
$ws.Range("N130").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = $null
$ws.Range("N131").Value = 0
$ws.Range("H132").Value = 1587.1613
$ws.Range("I132").Value = 1091.04
$ws.Range("K132").Value = 3273.12
$ws.Range("M132").Value = -743.1199999999999
$ws.Range("H133").Value = 38499.2
$ws.Range("J133").Value = 38499.2
$ws.Range("L133").Value = 38499.2
$ws.Range("N133").Value = -43559.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 34631.938
$ws.Range("I107").Value = 39079.355
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 39079.355
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = -37159.355
$ws.Range("N107").Value = -7340
$ws.Range("H134").Value = 2331.919
$ws.Range("I134").Value = 1946.4688
$ws.Range("K134").Value = 5839.4064
$ws.Range("M134").Value = -3304.4064

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 27375
$ws.Range("I23").Value = 34833.332
$ws.Range("J23").Value = 5000
$ws.Range("K23").Value = 34833.332
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = -34593.332
$ws.Range("N23").Value = -5480
$ws.Range("H27").Value = 27375
$ws.Range("I27").Value = 34833.332
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 34833.332
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -34641.332
$ws.Range("N27").Value = -5384
$ws.Range("H31").Value = 23812290
$ws.Range("J31").Value = 3512.5881
$ws.Range("L31").Value = 3512.5881
$ws.Range("N31").Value = -4102.5881
$ws.Range("H34").Value = 23812290
$ws.Range("J34").Value = 3512.5881
$ws.Range("L34").Value = 3512.5881
$ws.Range("N34").Value = -3916.5881
$ws.Range("H60").Value = 10107.5
$ws.Range("I60").Value = 5143.3335
$ws.Range("K60").Value = 5143.3335
$ws.Range("M60").Value = -4632.3335
$ws.Range("H62").Value = 252005
$ws.Range("I62").Value = 252005
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 252005
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -251381
$ws.Range("H65").Value = 252005
$ws.Range("I65").Value = 252005
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 1260025
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -1256905
$ws.Range("H120").Value = 21000
$ws.Range("J120").Value = 21000
$ws.Range("L120").Value = 21000
$ws.Range("N120").Value = -28258
$ws.Range("H121").Value = 22431.666
$ws.Range("J121").Value = 22431.666
$ws.Range("L121").Value = 22431.666
$ws.Range("N121").Value = -25051.666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1599
$ws.Range("I33").Value = 458
$ws.Range("J33").Value = 3500.6667
$ws.Range("K33").Value = 2748
$ws.Range("L33").Value = 21004.0002
$ws.Range("M33").Value = -2465
$ws.Range("N33").Value = -21570.0002
$ws.Range("H34").Value = 1300
$ws.Range("I34").Value = 100
$ws.Range("J34").Value = 1400
$ws.Range("K34").Value = 300
$ws.Range("L34").Value = 4200
$ws.Range("M34").Value = -216
$ws.Range("N34").Value = -4368
$ws.Range("H39").Value = 3669.0588
$ws.Range("J39").Value = 3669.0588
$ws.Range("L39").Value = 11007.1764
$ws.Range("N39").Value = -11595.1764
$ws.Range("H55").Value = 6155.037
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 6372.5386
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 19117.6158
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -19471.6158
$ws.Range("H69").Value = 928.5
$ws.Range("I69").Value = 450
$ws.Range("J69").Value = 1008.25
$ws.Range("K69").Value = 1350
$ws.Range("L69").Value = 3024.75
$ws.Range("M69").Value = -539
$ws.Range("N69").Value = -4646.75
$ws.Range("H72").Value = 928.5
$ws.Range("I72").Value = 450
$ws.Range("J72").Value = 1008.25
$ws.Range("K72").Value = 4050
$ws.Range("L72").Value = 9074.25
$ws.Range("M72").Value = 6
$ws.Range("N72").Value = -17186.25
$ws.Range("H131").Value = 11238335
$ws.Range("J131").Value = 12196453
$ws.Range("L131").Value = 36589359
$ws.Range("N131").Value = -36599439
$ws.Range("H137").Value = 30305806
$ws.Range("J137").Value = 37040270
$ws.Range("L137").Value = 111120810
$ws.Range("N137").Value = -111131010

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6896.304
$ws.Range("I70").Value = 6555.091
$ws.Range("K70").Value = 6555.091
$ws.Range("M70").Value = -6285.091
$ws.Range("H73").Value = 6896.304
$ws.Range("I73").Value = 6555.091
$ws.Range("K73").Value = 6555.091
$ws.Range("M73").Value = -5619.091
$ws.Range("H102").Value = 1775.9412
$ws.Range("I102").Value = 1795.8077
$ws.Range("J102").Value = 1711.375
$ws.Range("K102").Value = 1795.8077
$ws.Range("L102").Value = 1711.375
$ws.Range("M102").Value = -173.8077000000001
$ws.Range("N102").Value = -4955.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 3208.5
$ws.Range("I17").Value = 5408
$ws.Range("J17").Value = 1009
$ws.Range("K17").Value = 5408
$ws.Range("L17").Value = 1009
$ws.Range("M17").Value = -5238
$ws.Range("N17").Value = -1349
$ws.Range("H60").Value = 8000
$ws.Range("I60").Value = 8000
$ws.Range("K60").Value = 8000
$ws.Range("M60").Value = -7491
$ws.Range("H109").Value = 29000
$ws.Range("J109").Value = 29000
$ws.Range("L109").Value = 29000
$ws.Range("N109").Value = -31774
